# Apply "Natmi following Dr Hou advice" data update to LR-pairs sheet (Tgfb1-Itgb6)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tgfb1"
$ws.Cells.Item(2, 3).Value = "Itgb6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 96.320746
$ws.Cells.Item(2, 8).Value = 288.962238
$ws.Cells.Item(2, 9).Value = 0.3809824610908788
$ws.Cells.Item(2, 10).Value = 0.3809824610908788
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.111678
$ws.Cells.Item(2, 14).Value = 0.335034
$ws.Cells.Item(2, 15).Value = 0.01499007198665366
$ws.Cells.Item(2, 16).Value = 0.01499007198665366
$ws.Cells.Item(2, 17).Value = 10.756908271788
$ws.Cells.Item(2, 18).Value = 96.812174446092
$ws.Cells.Item(2, 19).Value = 0.005710954517404749
$ws.Cells.Item(2, 20).Value = 0.00571095451740475

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tgfb1"
$ws.Cells.Item(3, 3).Value = "Itgb6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 96.320746
$ws.Cells.Item(3, 8).Value = 288.962238
$ws.Cells.Item(3, 9).Value = 0.3809824610908788
$ws.Cells.Item(3, 10).Value = 0.3809824610908788
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.242708666666667
$ws.Cells.Item(3, 14).Value = 6.728126
$ws.Cells.Item(3, 15).Value = 0.3010294270888212
$ws.Cells.Item(3, 16).Value = 0.3010294270888212
$ws.Cells.Item(3, 17).Value = 216.0193718339987
$ws.Cells.Item(3, 18).Value = 1944.174346505988
$ws.Cells.Item(3, 19).Value = 0.1146869319930764
$ws.Cells.Item(3, 20).Value = 0.1146869319930764

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tgfb1"
$ws.Cells.Item(4, 3).Value = "Itgb6"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 96.320746
$ws.Cells.Item(4, 8).Value = 288.962238
$ws.Cells.Item(4, 9).Value = 0.3809824610908788
$ws.Cells.Item(4, 10).Value = 0.3809824610908788
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02011033333333333
$ws.Cells.Item(4, 14).Value = 0.060331
$ws.Cells.Item(4, 15).Value = 0.002699326137128775
$ws.Cells.Item(4, 16).Value = 0.002699326137128775
$ws.Cells.Item(4, 17).Value = 1.937042308975333
$ws.Cells.Item(4, 18).Value = 17.433380780778
$ws.Cells.Item(4, 19).Value = 0.001028395915010256
$ws.Cells.Item(4, 20).Value = 0.001028395915010256

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Tgfb1"
$ws.Cells.Item(5, 3).Value = "Itgb6"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 96.320746
$ws.Cells.Item(5, 8).Value = 288.962238
$ws.Cells.Item(5, 9).Value = 0.3809824610908788
$ws.Cells.Item(5, 10).Value = 0.3809824610908788
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.075634
$ws.Cells.Item(5, 14).Value = 15.226902
$ws.Cells.Item(5, 15).Value = 0.6812811747873964
$ws.Cells.Item(5, 16).Value = 0.6812811747873964
$ws.Cells.Item(5, 17).Value = 488.888853302964
$ws.Cells.Item(5, 18).Value = 4399.999679726676
$ws.Cells.Item(5, 19).Value = 0.2595561786653874
$ws.Cells.Item(5, 20).Value = 0.2595561786653874

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tgfb1"
$ws.Cells.Item(6, 3).Value = "Itgb6"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 18.46467533333333
$ws.Cells.Item(6, 8).Value = 55.394026
$ws.Cells.Item(6, 9).Value = 0.07303429161291354
$ws.Cells.Item(6, 10).Value = 0.07303429161291354
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.111678
$ws.Cells.Item(6, 14).Value = 0.335034
$ws.Cells.Item(6, 15).Value = 0.01499007198665366
$ws.Cells.Item(6, 16).Value = 0.01499007198665366
$ws.Cells.Item(6, 17).Value = 2.062098011876
$ws.Cells.Item(6, 18).Value = 18.558882106884
$ws.Cells.Item(6, 19).Value = 0.001094789288771829
$ws.Cells.Item(6, 20).Value = 0.001094789288771829

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tgfb1"
$ws.Cells.Item(7, 3).Value = "Itgb6"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 18.46467533333333
$ws.Cells.Item(7, 8).Value = 55.394026
$ws.Cells.Item(7, 9).Value = 0.07303429161291354
$ws.Cells.Item(7, 10).Value = 0.07303429161291354
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.242708666666667
$ws.Cells.Item(7, 14).Value = 6.728126
$ws.Cells.Item(7, 15).Value = 0.3010294270888212
$ws.Cells.Item(7, 16).Value = 0.3010294270888212
$ws.Cells.Item(7, 17).Value = 41.41088739725289
$ws.Cells.Item(7, 18).Value = 372.6979865752759
$ws.Cells.Item(7, 19).Value = 0.02198547096207326
$ws.Cells.Item(7, 20).Value = 0.02198547096207326

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Tgfb1"
$ws.Cells.Item(8, 3).Value = "Itgb6"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 18.46467533333333
$ws.Cells.Item(8, 8).Value = 55.394026
$ws.Cells.Item(8, 9).Value = 0.07303429161291354
$ws.Cells.Item(8, 10).Value = 0.07303429161291354
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.02011033333333333
$ws.Cells.Item(8, 14).Value = 0.060331
$ws.Cells.Item(8, 15).Value = 0.002699326137128775
$ws.Cells.Item(8, 16).Value = 0.002699326137128775
$ws.Cells.Item(8, 17).Value = 0.3713307758451111
$ws.Cells.Item(8, 18).Value = 3.341976982606
$ws.Cells.Item(8, 19).Value = 0.0001971433722574224
$ws.Cells.Item(8, 20).Value = 0.0001971433722574224

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Tgfb1"
$ws.Cells.Item(9, 3).Value = "Itgb6"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 18.46467533333333
$ws.Cells.Item(9, 8).Value = 55.394026
$ws.Cells.Item(9, 9).Value = 0.07303429161291354
$ws.Cells.Item(9, 10).Value = 0.07303429161291354
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 5.075634
$ws.Cells.Item(9, 14).Value = 15.226902
$ws.Cells.Item(9, 15).Value = 0.6812811747873964
$ws.Cells.Item(9, 16).Value = 0.6812811747873964
$ws.Cells.Item(9, 17).Value = 93.719933920828
$ws.Cells.Item(9, 18).Value = 843.4794052874519
$ws.Cells.Item(9, 19).Value = 0.04975688798981102
$ws.Cells.Item(9, 20).Value = 0.04975688798981102

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Tgfb1"
$ws.Cells.Item(10, 3).Value = "Itgb6"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 122.909391
$ws.Cells.Item(10, 8).Value = 368.728173
$ws.Cells.Item(10, 9).Value = 0.4861499128584522
$ws.Cells.Item(10, 10).Value = 0.4861499128584522
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.111678
$ws.Cells.Item(10, 14).Value = 0.335034
$ws.Cells.Item(10, 15).Value = 0.01499007198665366
$ws.Cells.Item(10, 16).Value = 0.01499007198665366
$ws.Cells.Item(10, 17).Value = 13.726274968098
$ws.Cells.Item(10, 18).Value = 123.536474712882
$ws.Cells.Item(10, 19).Value = 0.0072874221900536
$ws.Cells.Item(10, 20).Value = 0.007287422190053601

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Tgfb1"
$ws.Cells.Item(11, 3).Value = "Itgb6"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 122.909391
$ws.Cells.Item(11, 8).Value = 368.728173
$ws.Cells.Item(11, 9).Value = 0.4861499128584522
$ws.Cells.Item(11, 10).Value = 0.4861499128584522
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.242708666666667
$ws.Cells.Item(11, 14).Value = 6.728126
$ws.Cells.Item(11, 15).Value = 0.3010294270888212
$ws.Cells.Item(11, 16).Value = 0.3010294270888212
$ws.Cells.Item(11, 17).Value = 275.649956410422
$ws.Cells.Item(11, 18).Value = 2480.849607693798
$ws.Cells.Item(11, 19).Value = 0.1463454297470602
$ws.Cells.Item(11, 20).Value = 0.1463454297470602

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Tgfb1"
$ws.Cells.Item(12, 3).Value = "Itgb6"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 122.909391
$ws.Cells.Item(12, 8).Value = 368.728173
$ws.Cells.Item(12, 9).Value = 0.4861499128584522
$ws.Cells.Item(12, 10).Value = 0.4861499128584522
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.02011033333333333
$ws.Cells.Item(12, 14).Value = 0.060331
$ws.Cells.Item(12, 15).Value = 0.002699326137128775
$ws.Cells.Item(12, 16).Value = 0.002699326137128775
$ws.Cells.Item(12, 17).Value = 2.471748822807
$ws.Cells.Item(12, 18).Value = 22.245739405263
$ws.Cells.Item(12, 19).Value = 0.001312277166341696
$ws.Cells.Item(12, 20).Value = 0.001312277166341696

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Tgfb1"
$ws.Cells.Item(13, 3).Value = "Itgb6"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 122.909391
$ws.Cells.Item(13, 8).Value = 368.728173
$ws.Cells.Item(13, 9).Value = 0.4861499128584522
$ws.Cells.Item(13, 10).Value = 0.4861499128584522
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 5.075634
$ws.Cells.Item(13, 14).Value = 15.226902
$ws.Cells.Item(13, 15).Value = 0.6812811747873964
$ws.Cells.Item(13, 16).Value = 0.6812811747873964
$ws.Cells.Item(13, 17).Value = 623.8430838788939
$ws.Cells.Item(13, 18).Value = 5614.587754910045
$ws.Cells.Item(13, 19).Value = 0.3312047837549967
$ws.Cells.Item(13, 20).Value = 0.3312047837549967

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Tgfb1"
$ws.Cells.Item(14, 3).Value = "Itgb6"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 15.127183
$ws.Cells.Item(14, 8).Value = 45.381549
$ws.Cells.Item(14, 9).Value = 0.05983333443775553
$ws.Cells.Item(14, 10).Value = 0.05983333443775553
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.111678
$ws.Cells.Item(14, 14).Value = 0.335034
$ws.Cells.Item(14, 15).Value = 0.01499007198665366
$ws.Cells.Item(14, 16).Value = 0.01499007198665366
$ws.Cells.Item(14, 17).Value = 1.689373543074
$ws.Cells.Item(14, 18).Value = 15.204361887666
$ws.Cells.Item(14, 19).Value = 0.0008969059904234787
$ws.Cells.Item(14, 20).Value = 0.0008969059904234788

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Tgfb1"
$ws.Cells.Item(15, 3).Value = "Itgb6"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 15.127183
$ws.Cells.Item(15, 8).Value = 45.381549
$ws.Cells.Item(15, 9).Value = 0.05983333443775553
$ws.Cells.Item(15, 10).Value = 0.05983333443775553
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.242708666666667
$ws.Cells.Item(15, 14).Value = 6.728126
$ws.Cells.Item(15, 15).Value = 0.3010294270888212
$ws.Cells.Item(15, 16).Value = 0.3010294270888212
$ws.Cells.Item(15, 17).Value = 33.92586441635267
$ws.Cells.Item(15, 18).Value = 305.332779747174
$ws.Cells.Item(15, 19).Value = 0.01801159438661138
$ws.Cells.Item(15, 20).Value = 0.01801159438661138

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Tgfb1"
$ws.Cells.Item(16, 3).Value = "Itgb6"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 15.127183
$ws.Cells.Item(16, 8).Value = 45.381549
$ws.Cells.Item(16, 9).Value = 0.05983333443775553
$ws.Cells.Item(16, 10).Value = 0.05983333443775553
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.02011033333333333
$ws.Cells.Item(16, 14).Value = 0.060331
$ws.Cells.Item(16, 15).Value = 0.002699326137128775
$ws.Cells.Item(16, 16).Value = 0.002699326137128775
$ws.Cells.Item(16, 17).Value = 0.3042126925243334
$ws.Cells.Item(16, 18).Value = 2.737914232719
$ws.Cells.Item(16, 19).Value = 0.0001615096835194007
$ws.Cells.Item(16, 20).Value = 0.0001615096835194007

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Tgfb1"
$ws.Cells.Item(17, 3).Value = "Itgb6"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 15.127183
$ws.Cells.Item(17, 8).Value = 45.381549
$ws.Cells.Item(17, 9).Value = 0.05983333443775553
$ws.Cells.Item(17, 10).Value = 0.05983333443775553
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 5.075634
$ws.Cells.Item(17, 14).Value = 15.226902
$ws.Cells.Item(17, 15).Value = 0.6812811747873964
$ws.Cells.Item(17, 16).Value = 0.6812811747873964
$ws.Cells.Item(17, 17).Value = 76.780044359022
$ws.Cells.Item(17, 18).Value = 691.020399231198
$ws.Cells.Item(17, 19).Value = 0.04076332437720127
$ws.Cells.Item(17, 20).Value = 0.04076332437720127
